# KP-11725 D: Extension of questionnaire's translation files
# Adds a "Variable" column to both the "Translations" and "@@_question"
# sheets of the translation workbook, recording the question's variable
# name ("c1") alongside each existing translation row.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Translations")
$ws2 = $wb.Worksheets.Item("@@_question")

# --- Sheet "Translations": insert a new "Variable" column (B) ---
$ws1.Columns.Item(2).Insert()
$ws1.Columns.Item(2).ColumnWidth = 12.8
$ws1.Range("B1").Value = "Variable"
$ws1.Range("B2").Value = "c1"
$ws1.Range("B3").Value = "c1"
$ws1.Range("B4").Value = "c1"
$ws1.Range("B5").Value = "c1"

# --- Sheet "@@_question": insert a new "Variable" column (B) ---
$ws2.Columns.Item(2).Insert()
$ws2.Columns.Item(2).ColumnWidth = 6.8
$ws2.Range("B1").Value = "Variable"
$ws2.Range("B2").Value = "c1"

# --- Selections / active sheet (Translations becomes the active tab) ---
$ws2.Range("B2").Select()
$ws1.Activate()
$ws1.Range("B6").Select()
